$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$newSheet.Name = "Only age data"
$names = @(
    "Campylognathoides liasicus",
    "Campylognathoides zitteli",
    "Dorygnathus banthensis",
    "Eudimorphodon cromptonellus",
    "Eudimorphodon rosenfeldi",
    "Parapsicephalus purdoni",
    "Rhamphinion jenkinsi",
    "Dorygnathus purdoni",
    "Areripedctylus dehmi",
    "Caviramus schesaplanensis"
)
for ($i = 0; $i -lt $names.Count; $i++) {
    $cell = $newSheet.Cells.Item($i + 1, 1)
    $cell.Value = $names[$i]
}
Write-Host "ok"
